$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells F1:H1 - copy style (bold/centered/border) from existing header A1
$headers = @("KNN_Outliers_MAD", "SVM_Outliers_MAD", "RF_Outliers_MAD")
$cols = @(6, 7, 8)  # F, G, H

$ws.Range("A1").Copy()
for ($i = 0; $i -lt 3; $i++) {
    $col = $cols[$i]
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headers[$i]
    $cell.PasteSpecial(-4122)  # xlPasteFormats - reuse A1's header style
}

# New boolean FALSE cells F2:H5
for ($row = 2; $row -le 5; $row++) {
    for ($i = 0; $i -lt 3; $i++) {
        $col = $cols[$i]
        $ws.Cells.Item($row, $col).Value = $false
    }
}
